$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "JSU(-1.0355905243196528, 1.1752522890551038, 0.17114115529464008, 2.904084265298856)"
$ws.Range("C2").Value = "NIG(1.355267502764655, 1.0508519706661343, 3.702369343293296, 6.3092962029539645)"
$ws.Range("D2").Value = "NIG(2.7859255295559726, 1.9544332064097787, -0.5899427290163913, 5.4818348672087005)"
$ws.Range("E2").Value = "NCT(2.9014926361742397, 1.8848405886826147, -0.013655882714963697, 3.8764569044233665)"
